# Append a new daily price row (row 88) to the "Prices" sheet, duplicating
# the prior day's row (row 87) and only bumping the date in column A.
# Values in this sheet are stored as plain text (no numeric/date parsing),
# so we copy the existing text-typed row instead of assigning .Value
# directly (which would make Excel auto-coerce things like "5,384" into a
# number or "2025-05-28" into a date serial).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prices")

# Duplicate row 87 into row 88 (keeps every cell as plain text, no style
# changes) - this gives us B88:J88 already correct, matching B87:J87.
$ws.Range("A87:J87").Copy()
$ws.Range("A88:J88").PasteSpecial(-4104)

# Now fix up the date in A88. Writing the literal string via .Value would
# get reinterpreted as a date serial, so build it as a text formula result
# in a scratch cell and paste just the value back in - this keeps the cell
# a plain text string with no added number format / style.
$ws.Range("Z1").Formula = '="2025-05-28"'
$ws.Range("Z1").Copy()
$ws.Range("A88").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

$excel.CutCopyMode = 0
